# Apply weekly fruit/vegetable price update: rows 2-15 (excluding unchanged row 4)
# are re-populated from a reshuffled ordering of the original data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (now sources data originally in row 10)
$ws.Range("D2").Value = 44446
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 78000
$ws.Range("L2").Value = 78000
$ws.Range("M2").Value = 78000
$ws.Range("N2").Value = '$/caja 25 kilos'
$ws.Range("P2").Value = 3120
$ws.Range("Q2").Value = 25

# Row 3 (now sources data originally in row 11)
$ws.Range("D3").Value = 44446
$ws.Range("H3").Value = 'Inferno'
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 80000
$ws.Range("L3").Value = 80000
$ws.Range("M3").Value = 80000
$ws.Range("N3").Value = '$/caja 15 kilos'
$ws.Range("P3").Value = 5333
$ws.Range("Q3").Value = 15

# Row 5 (now sources data originally in row 14)
$ws.Range("D5").Value = 44425
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 75000
$ws.Range("L5").Value = 75000
$ws.Range("M5").Value = 75000
$ws.Range("P5").Value = 3000

# Row 6 (now sources data originally in row 3)
$ws.Range("D6").Value = 44343
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 36000
$ws.Range("L6").Value = 36000
$ws.Range("M6").Value = 36000
$ws.Range("P6").Value = 1440

# Row 7 (now sources data originally in row 13)
$ws.Range("D7").Value = 44421
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 75000
$ws.Range("L7").Value = 75000
$ws.Range("M7").Value = 75000
$ws.Range("P7").Value = 3000

# Row 8 (now sources data originally in row 6)
$ws.Range("D8").Value = 44460
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 95000
$ws.Range("L8").Value = 95000
$ws.Range("M8").Value = 95000
$ws.Range("N8").Value = '$/caja 25 kilos'
$ws.Range("P8").Value = 3800
$ws.Range("Q8").Value = 25

# Row 9 (now sources data originally in row 2)
$ws.Range("D9").Value = 44193
$ws.Range("K9").Value = 46000
$ws.Range("L9").Value = 46000
$ws.Range("M9").Value = 46000
$ws.Range("N9").Value = '$/caja 15 kilos'
$ws.Range("P9").Value = 3067
$ws.Range("Q9").Value = 15

# Row 10 (now sources data originally in row 7)
$ws.Range("D10").Value = 44449
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 80000
$ws.Range("L10").Value = 80000
$ws.Range("M10").Value = 80000
$ws.Range("P10").Value = 3200

# Row 11 (now sources data originally in row 8)
$ws.Range("D11").Value = 44449
$ws.Range("H11").Value = 'Americana (o)'
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 75000
$ws.Range("L11").Value = 75000
$ws.Range("M11").Value = 75000
$ws.Range("P11").Value = 5000

# Row 12 (now sources data originally in row 9)
$ws.Range("D12").Value = 44326
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("P12").Value = 1200

# Row 13 (now sources data originally in row 5)
$ws.Range("D13").Value = 44221
$ws.Range("J13").Value = 22
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 24545
$ws.Range("P13").Value = 982

# Row 14 (now sources data originally in row 15)
$ws.Range("D14").Value = 44340
$ws.Range("K14").Value = 35000
$ws.Range("L14").Value = 35000
$ws.Range("M14").Value = 35000
$ws.Range("P14").Value = 1400

# Row 15 (now sources data originally in row 12)
$ws.Range("D15").Value = 44474
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = 100000
$ws.Range("L15").Value = 100000
$ws.Range("M15").Value = 100000
$ws.Range("P15").Value = 4000

